$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview" (sheet1) - add row 3 for the new 61d0736d... file
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A3").Value = "61d0736d-fd40-4a60-bc52-4356430f6382ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$wsOverview.Range("B3").Value = "e2e\61d0736d-fd40-4a60-bc52-4356430f6382ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("D3").Value = ""
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-09-04 20:32:52"

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e265edb73a564f4e0b5f9fd0183a17415ee2ccad/e2e/61d0736d-fd40-4a60-bc52-4356430f6382ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md",
    "",
    "",
    "e2e\61d0736d-fd40-4a60-bc52-4356430f6382ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
) | Out-Null

$wsOverview.Range("B3").Style = "Hyperlink"
$wsOverview.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------
# Sheet "zh-cn" (sheet2) - add row 3
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$loZh.ListRows.Add() | Out-Null

$wsZh.Range("A3").Value = "61d0736d-fd40-4a60-bc52-4356430f6382ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("D3").Value = "e2e"
$wsZh.Range("E3").Value = "ht"
$wsZh.Range("F3").Value = "False"
$wsZh.Range("G3").Value = "61d0736d-fd40-4a60-bc52-4356430f6382oooooooooooooooooooooooooooooooooooooooo.9a371570a85f3da0a2781c1fac80d850cdd09720.zh-cn.xlf"
$wsZh.Range("H3").Value = "2016-09-04 20:32:47"
$wsZh.Range("I3").Value = ""
$wsZh.Range("J3").Value = ""
$wsZh.Range("K3").Value = "0001-01-01 00:00:00"
$wsZh.Range("L3").Value = ""
$wsZh.Range("M3").Value = "True"
$wsZh.Range("N3").Value = ""
$wsZh.Range("O3").Value = "False"
$wsZh.Range("P3").Value = ""

$wsZh.Hyperlinks.Add(
    $wsZh.Range("A3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e265edb73a564f4e0b5f9fd0183a17415ee2ccad/e2e/61d0736d-fd40-4a60-bc52-4356430f6382ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md",
    "",
    "",
    "61d0736d-fd40-4a60-bc52-4356430f6382ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
) | Out-Null

$wsZh.Range("A3").Style = "Hyperlink"
$wsZh.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------
# Sheet "de-de" (sheet3) - add row 3
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$loDe.ListRows.Add() | Out-Null

$wsDe.Range("A3").Value = "61d0736d-fd40-4a60-bc52-4356430f6382ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("D3").Value = "e2e"
$wsDe.Range("E3").Value = "ht"
$wsDe.Range("F3").Value = "False"
$wsDe.Range("G3").Value = "61d0736d-fd40-4a60-bc52-4356430f6382oooooooooooooooooooooooooooooooooooooooo.9a371570a85f3da0a2781c1fac80d850cdd09720.de-de.xlf"
$wsDe.Range("H3").Value = "2016-09-04 20:32:52"
$wsDe.Range("I3").Value = ""
$wsDe.Range("J3").Value = ""
$wsDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDe.Range("L3").Value = ""
$wsDe.Range("M3").Value = "True"
$wsDe.Range("N3").Value = ""
$wsDe.Range("O3").Value = "False"
$wsDe.Range("P3").Value = ""

$wsDe.Hyperlinks.Add(
    $wsDe.Range("A3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e265edb73a564f4e0b5f9fd0183a17415ee2ccad/e2e/61d0736d-fd40-4a60-bc52-4356430f6382ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md",
    "",
    "",
    "61d0736d-fd40-4a60-bc52-4356430f6382ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
) | Out-Null

$wsDe.Range("A3").Style = "Hyperlink"
$wsDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------
# Column width adjustments (Status columns got wider to fit
# "Ready for handoff")
# ---------------------------------------------------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 16.333333333333332
$wsOverview.Columns.Item(6).ColumnWidth = 16.333333333333332
$wsZh.Columns.Item(3).ColumnWidth = 16.333333333333332
$wsDe.Columns.Item(3).ColumnWidth = 16.333333333333332
